$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.461.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.190.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.86%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.558"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.189.35"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.51%  "
$ws.Range("E12").Value = "  +3.65%  "
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.713.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("E16").Value = "  +5.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.436.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.190.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "520.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.740"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.88%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +4.00%  "
$ws.Range("E28").Value = "  +2.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.08%  "
$ws.Range("E31").Value = "  +5.69%  "
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "511.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0427"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.50%  "
$ws.Range("E40").Value = "  +10.88%  "
$ws.Range("E41").Value = "  +2.05%  "
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("E43").Value = "  +7.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0678"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +15.28%  "
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.903.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +14.32%  "
$ws.Range("E49").Value = "  +3.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.32%  "
